# Trade #44 closed at 2026-02-17 13:27:53 - unknown UNKNOWN +0.000%
#
# Updates the "Summary" and "Strategy Status" roll-up sheets with the
# latest aggregate numbers, and appends the newly-closed MarketMaking
# trade (#44) as a new row on both the "All Trades" and "MarketMaking"
# trade-log sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Summary sheet roll-up numbers
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1197.57    # Current Capital
$summary.Range("B4").Value = -2.43      # Total P&L $
$summary.Range("B5").Value = -1.1       # Total P&L %
$summary.Range("B6").Value = 44         # Total Trades
$summary.Range("B7").Value = 18         # Winning Trades
$summary.Range("B9").Value = 40.91      # Win Rate %

# ---------------------------------------------------------------------
# 2. Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 97.56999999999999   # Capital
$status.Range("D4").Value = 44                  # Trades
$status.Range("E4").Value = -2.43               # P&L $
$status.Range("F4").Value = -2.43               # P&L %
$status.Range("G4").Value = 40.91               # Win Rate %

# ---------------------------------------------------------------------
# 3. Append the new closed trade (row 45) to both trade-log sheets
# ---------------------------------------------------------------------
function Add-ClosedTrade45($ws) {
    $ws.Cells.Item(45, 1).Value = 44

    # Date/time columns are stored as plain text in this workbook, so
    # force a text format before assigning - otherwise Excel would
    # auto-convert the ISO-looking strings into date/time serials.
    $cDate = $ws.Cells.Item(45, 2)
    $cDate.NumberFormat = "@"
    $cDate.Value = "2026-02-17"

    $cTime = $ws.Cells.Item(45, 3)
    $cTime.NumberFormat = "@"
    $cTime.Value = "13:27:46"

    $ws.Cells.Item(45, 4).Value = "MarketMaking"
    $ws.Cells.Item(45, 5).Value = "UP"
    $ws.Cells.Item(45, 6).Value = 0.83
    $ws.Cells.Item(45, 7).Value = 0.89
    $ws.Cells.Item(45, 8).Value = "CLOSED"
    $ws.Cells.Item(45, 9).Value = 7.2289
    $ws.Cells.Item(45, 10).Value = 0.06
    $ws.Cells.Item(45, 11).Value = 97.56999999999999
    $ws.Cells.Item(45, 12).Value = 0
    $ws.Cells.Item(45, 13).Value = 0
    $ws.Cells.Item(45, 14).Value = 0.6
    $ws.Cells.Item(45, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(45, 16).Value = "early_exit"
    $ws.Cells.Item(45, 17).Value = 0.1
}

Add-ClosedTrade45 $wb.Worksheets.Item("All Trades")
Add-ClosedTrade45 $wb.Worksheets.Item("MarketMaking")
